# Fruta / hortaliza, semanal
#
# The weekly refresh inserts one new daily-price record for
# "Macroferia Regional de Talca - Espárragos" as row 95, pushing the
# previously-last record (old row 101) down to become row 102.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 95; everything from old row 95 onward shifts
# down by one (old row 101 -> new row 102), matching the diff exactly.
$ws.Rows("95:95").Insert()

# Populate the newly inserted row 95 with the new record's data.
$ws.Cells.Item(95, 1).Value  = 5
$ws.Cells.Item(95, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(95, 3).Value  = "Maule"
$ws.Cells.Item(95, 4).Value  = 45180
$ws.Cells.Item(95, 5).Value  = 7
$ws.Cells.Item(95, 6).Value  = 300000000
$ws.Cells.Item(95, 7).Value  = "Espárragos"
$ws.Cells.Item(95, 8).Value  = "Sin especificar"
$ws.Cells.Item(95, 9).Value  = "Primera"
$ws.Cells.Item(95, 10).Value = 1000
$ws.Cells.Item(95, 11).Value = 2300
$ws.Cells.Item(95, 12).Value = 2300
$ws.Cells.Item(95, 13).Value = 2300
$ws.Cells.Item(95, 14).Value = "$/kilo"
$ws.Cells.Item(95, 15).Value = "Provincia de Linares"
$ws.Cells.Item(95, 16).Value = 2300
$ws.Cells.Item(95, 17).Value = 1
$ws.Cells.Item(95, 18).Value = "Hortaliza"
